$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "19+15="
$t.Cell(1, 2).Range.Text = "92-35="
$t.Cell(1, 3).Range.Text = "38+37="
$t.Cell(1, 4).Range.Text = "94-17="
$t.Cell(1, 5).Range.Text = "15+26="
$t.Cell(2, 1).Range.Text = "60-32="
$t.Cell(2, 2).Range.Text = "20-12="
$t.Cell(2, 3).Range.Text = "25-16="
$t.Cell(2, 4).Range.Text = "40-14="
$t.Cell(2, 5).Range.Text = "59+5="
$t.Cell(3, 1).Range.Text = "52-25="
$t.Cell(3, 2).Range.Text = "50-14="
$t.Cell(3, 3).Range.Text = "66-18="
$t.Cell(3, 4).Range.Text = "91-72="
$t.Cell(3, 5).Range.Text = "43+29="
$t.Cell(4, 1).Range.Text = "27+55="
$t.Cell(4, 2).Range.Text = "31-22="
$t.Cell(4, 3).Range.Text = "71-65="
$t.Cell(4, 4).Range.Text = "16+15="
$t.Cell(4, 5).Range.Text = "12+9="
$t.Cell(5, 1).Range.Text = "3+48="
$t.Cell(5, 2).Range.Text = "83-5="
$t.Cell(5, 3).Range.Text = "84-65="
$t.Cell(5, 4).Range.Text = "71-37="
$t.Cell(5, 5).Range.Text = "80-7="
$t.Cell(6, 1).Range.Text = "16+65="
$t.Cell(6, 2).Range.Text = "70-59="
$t.Cell(6, 3).Range.Text = "69+4="
$t.Cell(6, 4).Range.Text = "66+8="
$t.Cell(6, 5).Range.Text = "84-36="
$t.Cell(7, 1).Range.Text = "18+36="
$t.Cell(7, 2).Range.Text = "19+18="
$t.Cell(7, 3).Range.Text = "33-17="
$t.Cell(7, 4).Range.Text = "26+38="
$t.Cell(7, 5).Range.Text = "51-38="
$t.Cell(8, 1).Range.Text = "70-46="
$t.Cell(8, 2).Range.Text = "50-4="
$t.Cell(8, 3).Range.Text = "64-45="
$t.Cell(8, 4).Range.Text = "20-3="
$t.Cell(8, 5).Range.Text = "71-62="
$t.Cell(9, 1).Range.Text = "13+39="
$t.Cell(9, 2).Range.Text = "23-19="
$t.Cell(9, 3).Range.Text = "16+18="
$t.Cell(9, 4).Range.Text = "80-64="
$t.Cell(9, 5).Range.Text = "47-28="
$t.Cell(10, 1).Range.Text = "17+14="
$t.Cell(10, 2).Range.Text = "15+66="
$t.Cell(10, 3).Range.Text = "45+46="
$t.Cell(10, 4).Range.Text = "70-28="
$t.Cell(10, 5).Range.Text = "35+49="
$t.Cell(11, 1).Range.Text = "27-9="
$t.Cell(11, 2).Range.Text = "22+69="
$t.Cell(11, 3).Range.Text = "59+24="
$t.Cell(11, 4).Range.Text = "25+59="
$t.Cell(11, 5).Range.Text = "36+56="
$t.Cell(12, 1).Range.Text = "86+7="
$t.Cell(12, 2).Range.Text = "66-49="
$t.Cell(12, 3).Range.Text = "14+77="
$t.Cell(12, 4).Range.Text = "97-79="
$t.Cell(12, 5).Range.Text = "51-38="
$t.Cell(13, 1).Range.Text = "59+33="
$t.Cell(13, 2).Range.Text = "80-48="
$t.Cell(13, 3).Range.Text = "90-23="
$t.Cell(13, 4).Range.Text = "9+8="
$t.Cell(13, 5).Range.Text = "27+56="
$t.Cell(14, 1).Range.Text = "97-48="
$t.Cell(14, 2).Range.Text = "91-13="
$t.Cell(14, 3).Range.Text = "76+8="
$t.Cell(14, 4).Range.Text = "5+19="
$t.Cell(14, 5).Range.Text = "5+27="
$t.Cell(15, 1).Range.Text = "13+39="
$t.Cell(15, 2).Range.Text = "27+39="
$t.Cell(15, 3).Range.Text = "83-45="
$t.Cell(15, 4).Range.Text = "46+27="
$t.Cell(15, 5).Range.Text = "7+76="
$t.Cell(16, 1).Range.Text = "74-26="
$t.Cell(16, 2).Range.Text = "12+29="
$t.Cell(16, 3).Range.Text = "26+38="
$t.Cell(16, 4).Range.Text = "92-87="
$t.Cell(16, 5).Range.Text = "19+69="
$t.Cell(17, 1).Range.Text = "80-36="
$t.Cell(17, 2).Range.Text = "74+18="
$t.Cell(17, 3).Range.Text = "83-76="
$t.Cell(17, 4).Range.Text = "28+58="
$t.Cell(17, 5).Range.Text = "91-2="
$t.Cell(18, 1).Range.Text = "14+27="
$t.Cell(18, 2).Range.Text = "54-25="
$t.Cell(18, 3).Range.Text = "55-28="
$t.Cell(18, 4).Range.Text = "53-37="
$t.Cell(18, 5).Range.Text = "59+35="
$t.Cell(19, 1).Range.Text = "45+26="
$t.Cell(19, 2).Range.Text = "92-88="
$t.Cell(19, 3).Range.Text = "78+16="
$t.Cell(19, 4).Range.Text = "74-46="
$t.Cell(19, 5).Range.Text = "43-4="
$t.Cell(20, 1).Range.Text = "91-56="
$t.Cell(20, 2).Range.Text = "27+14="
$t.Cell(20, 3).Range.Text = "8+17="
$t.Cell(20, 4).Range.Text = "2+89="
$t.Cell(20, 5).Range.Text = "26+17="
